$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new value would otherwise be
# auto-detected as a number by Excel (e.g. "1.003"), to preserve the
# original inlineStr/text semantics of column D.
$textCells = @("D4","D5","D6","D7","D9","D10","D11","D13","D14","D15","D17","D18","D21","D22","D24","D26","D27","D28","D29","D30","D32","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D46","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "22.423.12"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "1.573.93"
$ws.Range("E3").Value = "  +0.41%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "1.003"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").Value = "291.26"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "0.3765"
$ws.Range("E7").Value = "  +2.69%  "
$ws.Range("E8").Value = "  +0.77%  "
$ws.Range("D9").Value = "0.3427"
$ws.Range("E9").Value = "  +1.60%  "
$ws.Range("D10").Value = "1.162"
$ws.Range("E10").Value = "  -0.51%  "
$ws.Range("D11").Value = "0.07674"
$ws.Range("E11").Value = "  +1.40%  "
$ws.Range("E12").Value = "  +0.21%  "
$ws.Range("D13").Value = "21.28"
$ws.Range("E13").Value = "  +1.08%  "
$ws.Range("D14").Value = "6.008"
$ws.Range("E14").Value = "  -0.63%  "
$ws.Range("D15").Value = "6.944"
$ws.Range("E15").Value = "  +1.33%  "
$ws.Range("D16").Value = "1.574.90"
$ws.Range("E16").Value = "  +0.16%  "
$ws.Range("D17").Value = "0.00001136"
$ws.Range("D18").Value = "90.34"
$ws.Range("E18").Value = "  +1.41%  "
$ws.Range("E19").Value = "  +0.66%  "
$ws.Range("D21").Value = "16.80"
$ws.Range("E21").Value = "  +2.59%  "
$ws.Range("D22").Value = "6.232"
$ws.Range("E22").Value = "  -0.27%  "
$ws.Range("E23").Value = "  +0.76%  "
$ws.Range("D24").Value = "2.429"
$ws.Range("E24").Value = "  +1.98%  "
$ws.Range("D25").Value = "22.426.97"
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("D26").Value = "2.744"
$ws.Range("E26").Value = "  -6.14%  "
$ws.Range("D27").Value = "20.32"
$ws.Range("E27").Value = "  +2.18%  "
$ws.Range("D28").Value = "146.35"
$ws.Range("E28").Value = "  +0.25%  "
$ws.Range("D29").Value = "5.023"
$ws.Range("E29").Value = "  +1.43%  "
$ws.Range("D30").Value = "126.39"
$ws.Range("E30").Value = "  +1.22%  "
$ws.Range("D31").Value = "1.750.47"
$ws.Range("E31").Value = "  +0.24%  "
$ws.Range("D32").Value = "6.214"
$ws.Range("E32").Value = "  -0.41%  "
$ws.Range("E33").Value = "  +1.81%  "
$ws.Range("D34").Value = "1.001"
$ws.Range("E34").Value = "  +1.84%  "
$ws.Range("D35").Value = "10.06"
$ws.Range("E35").Value = "  -2.50%  "
$ws.Range("D36").Value = "0.08578"
$ws.Range("E36").Value = "  +1.98%  "
$ws.Range("D37").Value = "0.02555"
$ws.Range("E37").Value = "  +1.12%  "
$ws.Range("D38").Value = "0.2318"
$ws.Range("E38").Value = "  +1.02%  "
$ws.Range("D39").Value = "0.06581"
$ws.Range("E39").Value = "  +1.37%  "
$ws.Range("D40").Value = "1.330"
$ws.Range("E40").Value = "  +7.14%  "
$ws.Range("D41").Value = "5.476"
$ws.Range("E41").Value = "  -0.61%  "
$ws.Range("D42").Value = "11.58"
$ws.Range("E42").Value = "  -1.48%  "
$ws.Range("D43").Value = "0.6458"
$ws.Range("E43").Value = "  +1.28%  "
$ws.Range("D44").Value = "14.14"
$ws.Range("E44").Value = "  -1.95%  "
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "0.6024"
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "3.800"
$ws.Range("E47").Value = "  +0.72%  "
$ws.Range("D48").Value = "1.297"
$ws.Range("E48").Value = "  +9.57%  "
$ws.Range("D49").Value = "2.086"
$ws.Range("E49").Value = "  -1.32%  "
$ws.Range("D50").Value = "125.46"
$ws.Range("E50").Value = "  +3.66%  "
$ws.Range("D51").Value = "0.07327"
$ws.Range("E51").Value = "  +0.85%  "
